# Sync BaseData sheet with latest measurement campaign (LM):
#  - Existing rows 2 & 3 (previously Denges 2011 / Denges 2012) become the
#    Ceneri site for years 2018 and 2017.
#  - Two new rows are appended re-using the original Denges 2017 / 2018
#    values (copied from the existing row layout/format).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BaseData")

# Copy the formatting/values of rows 2:3 down into rows 4:5 before touching
# the B/C columns, so the new rows inherit the same style pattern as the
# existing data rows.
$ws.Range("A2:N3").Copy($ws.Range("A4")) | Out-Null
$excel.CutCopyMode = $false

# Update row 2: Denges/2011 -> Ceneri/2018
$ws.Range("B2").Value = "Ceneri"
$ws.Range("C2").Value = 2018

# Update row 3: Denges/2012 -> Ceneri/2017
$ws.Range("B3").Value = "Ceneri"
$ws.Range("C3").Value = 2017

# New row 4: Denges/2017
$ws.Range("B4").Value = "Denges"
$ws.Range("C4").Value = 2017

# New row 5: Denges/2018
$ws.Range("B5").Value = "Denges"
$ws.Range("C5").Value = 2018

$ws.Range("K6").Select() | Out-Null
